# Update the Training Dashboard progress-as-of date from 03-Nov-2025 to
# 04-Nov-2025, and decrement the corresponding "PERIOD TO EXPIRE" day
# counts by 1 (H3: 138 -> 137, H4: 137 -> 136).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$ws.Range("H3").Value = 137
$ws.Range("I3").Value = "'04-Nov-2025"

$ws.Range("H4").Value = 136
$ws.Range("I4").Value = "'04-Nov-2025"
